# Update grading numbers for the GradedExercise workbook, plus a comment in code.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Functionality section (rows 4-6) --
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 12
# G6/H6 are formulas (SUM/MIN) and recalculate automatically.

# -- Engineering section (rows 8-13) --
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 6
# G13/H13 are formulas and recalculate automatically.

# -- User Input section (rows 15-22) --
$ws.Range("G16").Value = 1
# G22/H22 are formulas and recalculate automatically.

# -- Extra section (rows 24-29): clear the per-criteria scores that no
# -- longer apply, keep only the last one (now worth 2 instead of 5).
$ws.Range("G24").ClearContents()
$ws.Range("G25").ClearContents()
$ws.Range("G26").ClearContents()
$ws.Range("G27").ClearContents()
$ws.Range("G28").Value = 2
# G29/H29 are formulas and recalculate automatically.

# -- Totals (row 32) and grade (row 34) recalculate automatically via
# -- their existing SUM/MIN/TOTAL_POINTS formulas.

# Move the active selection to reflect where grading left off.
$ws.Range("H28").Select()
